$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 364.07144
$ws.Range("I33").Value = 305.42856
$ws.Range("K33").Value = 305.42856
$ws.Range("M33").Value = -76.42856

$ws.Range("H62").Value = 7703.1055
$ws.Range("I62").Value = 6842.143
$ws.Range("K62").Value = 6842.143
$ws.Range("M62").Value = -6218.143

$ws.Range("H65").Value = 7703.1055
$ws.Range("I65").Value = 6842.143
$ws.Range("K65").Value = 34210.715
$ws.Range("M65").Value = -31090.715

$ws.Range("H80").Value = 6642.85
$ws.Range("I80").Value = 885.7
$ws.Range("K80").Value = 2657.1
$ws.Range("M80").Value = -1659.1

$ws.Range("H83").Value = 6642.85
$ws.Range("I83").Value = 885.7
$ws.Range("K83").Value = 7971.3
$ws.Range("M83").Value = -2979.3

$ws.Range("H107").Value = 1376.579
$ws.Range("I107").Value = 940.4286
$ws.Range("J107").Value = 2597.8
$ws.Range("K107").Value = 940.4286
$ws.Range("L107").Value = 2597.8
$ws.Range("M107").Value = 979.5714
$ws.Range("N107").Value = -6437.8

$ws.Range("H118").Value = 1888.3125
$ws.Range("I118").Value = 571.4
$ws.Range("J118").Value = 4083.1667
$ws.Range("K118").Value = 1714.2
$ws.Range("L118").Value = 12249.5001
$ws.Range("M118").Value = -57.19999999999982
$ws.Range("N118").Value = -15563.5001

$ws.Range("H137").Value = 45064
$ws.Range("I137").Value = 2890
$ws.Range("K137").Value = 8670
$ws.Range("M137").Value = -6120

$ws.Range("H141").Value = 1837.1428
$ws.Range("I141").Value = 1785.0834
$ws.Range("K141").Value = 5355.2502
$ws.Range("M141").Value = -175.2502000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3844
$ws.Range("I45").Value = 2632.8
$ws.Range("J45").Value = 9900
$ws.Range("K45").Value = 2632.8
$ws.Range("L45").Value = 9900
$ws.Range("M45").Value = -2255.8
$ws.Range("N45").Value = -10654

$ws.Range("H110").Value = 23950.05
$ws.Range("I110").Value = 25374.916
$ws.Range("K110").Value = 25374.916
$ws.Range("M110").Value = -23329.916

$ws.Range("H122").Value = 2065.2
$ws.Range("I122").Value = 1899.5714
$ws.Range("J122").Value = 2451.6667
$ws.Range("K122").Value = 5698.7142
$ws.Range("L122").Value = 7355.000100000001
$ws.Range("M122").Value = -3248.7142
$ws.Range("N122").Value = -12255.0001

$ws.Range("H127").Value = 96375
$ws.Range("J127").Value = 96375
$ws.Range("L127").Value = 96375
$ws.Range("N127").Value = -106295

$ws.Range("H128").Value = 90429
$ws.Range("J128").Value = 90429
$ws.Range("L128").Value = 90429
$ws.Range("N128").Value = -100389

$ws.Range("H131").Value = 103578.8
$ws.Range("J131").Value = 103578.8
$ws.Range("L131").Value = 103578.8
$ws.Range("N131").Value = -113658.8

$ws.Range("H132").Value = 1537.2285
$ws.Range("I132").Value = 1335.8572
$ws.Range("K132").Value = 4007.5716
$ws.Range("M132").Value = -1477.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9774.9
$ws.Range("I20").Value = 13104.571
$ws.Range("K20").Value = 13104.571
$ws.Range("M20").Value = -12857.571

$ws.Range("H99").Value = 779.3333
$ws.Range("I99").Value = 766.41174
$ws.Range("K99").Value = 766.41174
$ws.Range("M99").Value = 731.58826

$ws.Range("H105").Value = 3426.6
$ws.Range("I105").Value = 3069.375
$ws.Range("K105").Value = 3069.375
$ws.Range("M105").Value = -1322.375

$ws.Range("H107").Value = 2355.8572
$ws.Range("I107").Value = 2531.76
$ws.Range("J107").Value = 890
$ws.Range("K107").Value = 2531.76
$ws.Range("L107").Value = 890
$ws.Range("M107").Value = -611.7600000000002
$ws.Range("N107").Value = -4730

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4167578
$ws.Range("I31").Value = 5000888
$ws.Range("K31").Value = 5000888
$ws.Range("M31").Value = -5000593

$ws.Range("H34").Value = 4167578
$ws.Range("I34").Value = 5000888
$ws.Range("K34").Value = 5000888
$ws.Range("M34").Value = -5000686

$ws.Range("H93").Value = 31135.666
$ws.Range("I93").Value = 31135.666
$ws.Range("K93").Value = 31135.666
$ws.Range("M93").Value = -29263.666

$ws.Range("H99").Value = 7687.222
$ws.Range("I99").Value = 6546.3335
$ws.Range("J99").Value = 9969
$ws.Range("K99").Value = 6546.3335
$ws.Range("L99").Value = 9969
$ws.Range("M99").Value = -5048.3335
$ws.Range("N99").Value = -12965

$ws.Range("H122").Value = 1331.7693
$ws.Range("I122").Value = 1331.7693
$ws.Range("K122").Value = 3995.3079
$ws.Range("M122").Value = -1545.3079

$ws.Range("H126").Value = 7687.222
$ws.Range("I126").Value = 6546.3335
$ws.Range("J126").Value = 9969
$ws.Range("K126").Value = 19639.0005
$ws.Range("L126").Value = 29907
$ws.Range("M126").Value = -17169.0005
$ws.Range("N126").Value = -34847

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 206.78947
$ws.Range("I12").Value = 206.16667
$ws.Range("J12").Value = 207.07692
$ws.Range("K12").Value = 618.50001
$ws.Range("L12").Value = 621.23076
$ws.Range("M12").Value = -445.50001
$ws.Range("N12").Value = -967.23076

$ws.Range("H18").Value = 393.6
$ws.Range("I18").Value = 199.5
$ws.Range("J18").Value = 523
$ws.Range("K18").Value = 598.5
$ws.Range("L18").Value = 1569
$ws.Range("M18").Value = -429.5
$ws.Range("N18").Value = -1907

$ws.Range("H33").Value = 425
$ws.Range("I33").Value = 200.875
$ws.Range("J33").Value = 1022.6667
$ws.Range("K33").Value = 1205.25
$ws.Range("L33").Value = 6136.0002
$ws.Range("M33").Value = -922.25
$ws.Range("N33").Value = -6702.0002

$ws.Range("H113").Value = 1108.5
$ws.Range("J113").Value = 1213.4286
$ws.Range("L113").Value = 3640.2858
$ws.Range("N113").Value = -7980.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9998
$ws.Range("I70").Value = 9998
$ws.Range("K70").Value = 9998
$ws.Range("M70").Value = -9728

$ws.Range("H73").Value = 9998
$ws.Range("I73").Value = 9998
$ws.Range("K73").Value = 9998
$ws.Range("M73").Value = -9062

$ws.Range("H80").Value = 10731.462
$ws.Range("I80").Value = 5399.5
$ws.Range("J80").Value = 13101.223
$ws.Range("K80").Value = 5399.5
$ws.Range("L80").Value = 13101.223
$ws.Range("M80").Value = -4401.5
$ws.Range("N80").Value = -15097.223

$ws.Range("H83").Value = 10731.462
$ws.Range("I83").Value = 5399.5
$ws.Range("J83").Value = 13101.223
$ws.Range("K83").Value = 26997.5
$ws.Range("L83").Value = 65506.115
$ws.Range("M83").Value = -22005.5
$ws.Range("N83").Value = -75490.11499999999

$ws.Range("H97").Value = 1675.6154
$ws.Range("I97").Value = 1947.5
$ws.Range("K97").Value = 1947.5
$ws.Range("M97").Value = -1451.5

$ws.Range("H102").Value = 18422.768
$ws.Range("I102").Value = 21882.916
$ws.Range("J102").Value = 4582.1665
$ws.Range("K102").Value = 21882.916
$ws.Range("L102").Value = 4582.1665
$ws.Range("M102").Value = -20260.916
$ws.Range("N102").Value = -7826.1665

$ws.Range("H113").Value = 1857.1428
$ws.Range("I113").Value = 1666.6666
$ws.Range("K113").Value = 1666.6666
$ws.Range("M113").Value = 503.3334

$ws.Range("H124").Value = 60000
$ws.Range("J124").Value = 60000
$ws.Range("L124").Value = 60000
$ws.Range("N124").Value = -69820

$ws.Range("H126").Value = 2314.2104
$ws.Range("I126").Value = 1839.3
$ws.Range("K126").Value = 5517.9
$ws.Range("M126").Value = -3047.9

$ws.Range("H132").Value = 3485.5
$ws.Range("I132").Value = 2924.5715
$ws.Range("J132").Value = 4467.125
$ws.Range("K132").Value = 8773.7145
$ws.Range("L132").Value = 13401.375
$ws.Range("M132").Value = -6243.7145
$ws.Range("N132").Value = -18461.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2415.6667
$ws.Range("I55").Value = 1483.4
$ws.Range("K55").Value = 1483.4
$ws.Range("M55").Value = -1310.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 25422.023
$ws.Range("I122").Value = 29158.553
$ws.Range("K122").Value = 87475.659
$ws.Range("M122").Value = -85025.659
